# Add "Link to Coverage" entries in column I and refresh the sheet view /
# row-height / column-width tweaks that came along with the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CV32E40X Counters")

# --- New "Link to Coverage" text for the HPM-counter related requirement rows ---
$ws.Range("I2").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.counters_cg.x_check_mcycle"
$ws.Range("I3").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.counters_cg.x_check_minstret"
$ws.Range("I4").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.cg_idx_wrapper_*.mhpm_cg.x_check_mhpm"
$ws.Range("I5").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.cg_idx_wrapper_*.inhibit_mix_cg.x_check_*"
$ws.Range("I6").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.counters_cg.x_check_mcycle"
$ws.Range("I8").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.counters_cg.x_check_minstret"
$ws.Range("I32").Value = "uvm_pkg.uvm_test_top.env.cov_model.counters_covg.counters_cg.x_minstret_in_dbg"

# --- Row heights ---
$ws.Rows(1).RowHeight = 30
$ws.Rows(2).RowHeight = 45
$ws.Rows(3).RowHeight = 45

# --- Column widths ---
$ws.Columns("B").ColumnWidth = 16.1667
$ws.Columns("H").ColumnWidth = 11.5
$ws.Columns("I").ColumnWidth = 38.5

# --- Sheet view / pane / selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("I34").Select()
